# Add data for 2022-04-26
# Rename sheet to reflect new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-04-18"

# Update the header label for the current-year column
$ws.Range("I1").Value = "2022 (through 04-18)"

# Update April (row 5) value for the current-year column
$ws.Range("I5").Value = 80

# Update the Total (row 14) for the current-year column
$ws.Range("I14").Value = 515
